$wb = $excel.ActiveWorkbook

# "All Published Values" sheet: append the newly captured BOC USD publish as row 3.
$wsValues = $wb.Worksheets.Item("All Published Values")

# Format the new row's cells as Text first so date/number-looking strings
# (e.g. "2026-01-02", "697.85") are stored verbatim instead of being
# auto-coerced into date serials / numbers, matching the other rows.
$wsValues.Range("A3:J3").NumberFormat = "@"

$wsValues.Range("A3").Value = "2026-01-02"
$wsValues.Range("B3").Value = "2026-01-02 18:35:17"
$wsValues.Range("C3").Value = "697.85"
$wsValues.Range("D3").Value = "697.85"
$wsValues.Range("E3").Value = "700.79"
$wsValues.Range("F3").Value = "700.79"
$wsValues.Range("G3").Value = "702.88"
$wsValues.Range("H3").Value = "2026/01/02 18:35:17"
$wsValues.Range("I3").Value = "2026-01-02 10:37:44"
$wsValues.Range("J3").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# "Daily Summary" sheet: bump the publishes count for 2026-01-02 from 1 to 2,
# now that a second rate was captured for that date.
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("B3").NumberFormat = "@"
$wsSummary.Range("B3").Value = "2"
